# Commit: update proforma in/out sheet title to "Results Summary and Inputs"
#
# The "Inputs and Outputs" worksheet's A1 banner cell changes text from
# "Inputs" to "Results Summary and Inputs". Sheet/tab names are untouched;
# only the displayed cell value changes (Excel reshuffles the shared-string
# table and cursor state as a natural side effect of the edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs and Outputs")

# Make sure we're editing/selecting on the right sheet, then update the title.
[void]$ws.Activate()
$ws.Range("A1").Value = "Results Summary and Inputs"

# Leave the cursor sitting on the cell we just edited (matches the resulting
# document's default/A1 cursor state instead of the stale prior selection).
[void]$ws.Range("A1").Select()
